# Update header for column C, add new column D "TotGoedningabDyr_kt_år_udbr"
# (udbredelse af gødning) with the potential/spread-out manure amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header in C1 to reflect it now being the "potential" total,
# and add the new header in D1.
$ws.Range("C1").Value = "TotGoednabDyr_kt_år_pot"
$ws.Range("D1").Value = "TotGoedningabDyr_kt_år_udbr"

$data = @(
    @{ Row = 2;  D = 3141.89058322933 },
    @{ Row = 3;  D = 4827.9579547877 },
    @{ Row = 4;  D = 0 },
    @{ Row = 5;  D = 0 },
    @{ Row = 6;  D = 357.1165508235 },
    @{ Row = 7;  D = 451.114122566825 },
    @{ Row = 8;  D = 978.944806013418 },
    @{ Row = 9;  D = 3456.07964155226 },
    @{ Row = 10; D = 0 },
    @{ Row = 11; D = 0 },
    @{ Row = 12; D = 1085.95563365066 },
    @{ Row = 13; D = 0 },
    @{ Row = 14; D = 0 },
    @{ Row = 15; D = 1085.95563365066 },
    @{ Row = 16; D = 0 },
    @{ Row = 17; D = 3866.49625958574 },
    @{ Row = 18; D = 2263.6160885719 },
    @{ Row = 19; D = 1085.95563365066 },
    @{ Row = 20; D = 0 },
    @{ Row = 21; D = 1085.95563365066 },
    @{ Row = 22; D = 0 },
    @{ Row = 23; D = 1085.95563365066 },
    @{ Row = 24; D = 0 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
